# Updated cryptos list on Tue Nov 12 14:43:11 UTC 2024 with GitHub Actions
#
# Refresh each cryptocurrency row's Price (column D) and Volume(1h) (column E)
# with newly-scraped figures. A handful of rows additionally had their
# ranking order reshuffled (coin name in column B and its coinranking.com
# link in column C moved to a different row), so those rows also update
# columns B and C to reflect the coin that now occupies that rank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '86.083.17'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +4.45%  '
# Row 3: Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.270.50'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +2.87%  '
# Row 4: TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.997'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.29%  '
# Row 5: Solana
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.64'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -3.83%  '
# Row 6: BNB
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '624.90'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.89%  '
# Row 7: Dogecoin
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.366'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +26.14%  '
# Row 8: XRP
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.651'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +11.94%  '
# Row 9: USDC
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.09%  '
# Row 10: LidoStakedEther
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.261.78'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +2.79%  '
# Row 11: Cardano
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.575'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -3.07%  '
# Row 12: TRON
$ws.Range("E12").Value = '  +6.85%  '
# Row 13: ShibaInu
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000255'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.55%  '
# Row 14: Avalanche
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.88'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +5.39%  '
# Row 15: WrappedliquidstakedEther2.0
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.852.50'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +2.41%  '
# Row 16: Toncoin
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.26'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.89%  '
# Row 17: WrappedBTC
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '85.842.67'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +4.53%  '
# Row 18: WrappedEther
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.258.46'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +2.75%  '
# Row 19: Chainlink
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.01'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.21%  '
# Row 20: SuiNetwork
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.99'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -8.05%  '
# Row 21: BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '427.69'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.45%  '
# Row 22: Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.91'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.05%  '
# Row 23: Polkadot
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.31'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +3.47%  '
# Row 24: LEO
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.18'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.57%  '
# Row 25: Aptos
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.37'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +3.50%  '
# Row 26: NEARProtocol
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.08'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -3.26%  '
# Row 27: WrappedeETH
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.437.14'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +2.78%  '
# Row 28: Litecoin
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '75.46'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.76%  '
# Row 29: PEPE
$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000128'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +5.66%  '
# Row 30: Dai
$ws.Range("B30").Value = 'Dai'
$ws.Range("C30").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.38%  '
# Row 31: Cronos
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.173'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +18.34%  '
# Row 32: Binance-PegBSC-USD
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.24%  '
# Row 33: InternetComputer(DFINITY)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.78'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.81%  '
# Row 34: Bittensor
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '545.79'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -4.44%  '
# Row 35: Fetch.AI
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.42'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -4.40%  '
# Row 36: PancakeSwap
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.94'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -2.16%  '
# Row 37: RenderToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.84'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +10.11%  '
# Row 38: Kaspa
$ws.Range("E38").Value = '  -10.00%  '
# Row 39: EthereumClassic
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '22.36'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.16%  '
# Row 40: FirstDigitalUSD
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.996'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.25%  '
# Row 41: WhiteBITCoin
$ws.Range("E41").Value = '  +3.54%  '
# Row 42: PolygonEcosystemToken
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.391'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -3.25%  '
# Row 43: Stacks
$ws.Range("E43").Value = '  -2.05%  '
# Row 44: dogwifhat
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.93'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.33%  '
# Row 45: Monero
$ws.Range("B45").Value = 'Monero'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '157.93'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.85%  '
# Row 46: USDe
$ws.Range("B46").Value = 'USDe'
$ws.Range("C46").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.14%  '
# Row 47: Aave
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '178.21'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -4.34%  '
# Row 48: OKB
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '44.01'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.49%  '
# Row 49: ImmutableX
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.29'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.66%  '
# Row 50: Filecoin
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.23'
$ws.Range("D50").ClearFormats()
# Row 51: ARBITRUM
$ws.Range("B51").Value = 'ARBITRUM'
$ws.Range("C51").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.619'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.39%  '
